$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "Created"
$ws.Range("G2").Value = "Yes"
